$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.052.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.67%  "

$ws.Range("D3").Value = "'1.665.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'216.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").Value = "'0.5107"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.2635"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "'0.06418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.46%  "

$ws.Range("D10").Value = "'21.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("D11").Value = "'0.07420"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "

$ws.Range("D12").Value = "'1.670.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.30%  "

$ws.Range("D13").Value = "'4.510"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").Value = "'0.5803"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").Value = "'0.000008542"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.68%  "

$ws.Range("D16").Value = "'64.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "'26.136.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").Value = "'4.915"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").Value = "'189.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.50%  "

$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'145.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").Value = "'7.622"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").Value = "'0.1206"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.76%  "

$ws.Range("E27").Value = "  +1.55%  "

$ws.Range("D28").Value = "'0.06405"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.69%  "

$ws.Range("D29").Value = "'1.296"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("D30").Value = "'1.314"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.72%  "

$ws.Range("D31").Value = "'3.520"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("D32").Value = "'3.506"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").Value = "'1.630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("D34").Value = "'1.015"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("D35").Value = "'0.6081"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.55%  "

$ws.Range("D36").Value = "'2.362"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").Value = "'2.647"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").Value = "'6.176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.37%  "

$ws.Range("D39").Value = "'0.01607"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.13%  "

$ws.Range("D40").Value = "'1.077.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("D43").Value = "'100.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.52%  "

$ws.Range("D44").Value = "'1.814.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "

$ws.Range("D45").Value = "'0.00000000114"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.15%  "

$ws.Range("D46").Value = "'56.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").Value = "'1.009"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").Value = "'8.084"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("D49").Value = "'0.05204"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").Value = "'0.4289"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").Value = "'5.948"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.49%  "
